# Weekly update: insert two new daily price records for
# Fruta, Terminal Hortofrutícola Agro Chillán - Kiwi
# New records are inserted as rows 95-96, pushing existing data down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 95 (shifts rows 95:198 down to 97:200)
$ws.Rows("95:96").Insert()

# ---- New row 95: Kiwi Hayward, Primera ----
$ws.Cells.Item(95, 1).Value = 7
$ws.Cells.Item(95, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(95, 3).Value = "Ñuble"
$ws.Cells.Item(95, 4).Value = 44790
$ws.Cells.Item(95, 5).Value = 16
$ws.Cells.Item(95, 6).Value = "Fruta"
$ws.Cells.Item(95, 7).Value = 100101
$ws.Cells.Item(95, 8).Value = "Berries"
$ws.Cells.Item(95, 9).Value = 100101007
$ws.Cells.Item(95, 10).Value = "Kiwi"
$ws.Cells.Item(95, 11).Value = "Hayward"
$ws.Cells.Item(95, 12).Value = "Primera"
$ws.Cells.Item(95, 13).Value = 80
$ws.Cells.Item(95, 14).Value = 7500
$ws.Cells.Item(95, 15).Value = 8000
$ws.Cells.Item(95, 16).Value = 7750
$ws.Cells.Item(95, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(95, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(95, 19).Value = 431
$ws.Cells.Item(95, 20).Value = 18

# ---- New row 96: Kiwi Hayward, Segunda ----
$ws.Cells.Item(96, 1).Value = 7
$ws.Cells.Item(96, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(96, 3).Value = "Ñuble"
$ws.Cells.Item(96, 4).Value = 44790
$ws.Cells.Item(96, 5).Value = 16
$ws.Cells.Item(96, 6).Value = "Fruta"
$ws.Cells.Item(96, 7).Value = 100101
$ws.Cells.Item(96, 8).Value = "Berries"
$ws.Cells.Item(96, 9).Value = 100101007
$ws.Cells.Item(96, 10).Value = "Kiwi"
$ws.Cells.Item(96, 11).Value = "Hayward"
$ws.Cells.Item(96, 12).Value = "Segunda"
$ws.Cells.Item(96, 13).Value = 120
$ws.Cells.Item(96, 14).Value = 6500
$ws.Cells.Item(96, 15).Value = 7000
$ws.Cells.Item(96, 16).Value = 6750
$ws.Cells.Item(96, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(96, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(96, 19).Value = 375
$ws.Cells.Item(96, 20).Value = 18
